$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Sheet 1 new rows
# Row 231
$ws1.Cells.Item(231, 1).Value = 45439.55208333334
$ws1.Cells.Item(231, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(231, 2).Value = "12-06-2024 09:15:00"
$ws1.Cells.Item(231, 3).Value = "hour"
$ws1.Cells.Item(231, 4).Value = "TATASTEEL.NS"
$ws1.Cells.Item(231, 5).Value = 45408.38541666666
$ws1.Cells.Item(231, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(231, 6).Value = 170.6999969482422
$ws1.Cells.Item(231, 7).Value = 45434.38541666666
$ws1.Cells.Item(231, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(231, 8).Value = 175.4499969482422
$ws1.Cells.Item(231, 9).Value = 45436.38541666666
$ws1.Cells.Item(231, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(231, 10).Value = 177.5500030517578
$ws1.Cells.Item(231, 11).Value = "High"
$ws1.Cells.Item(231, 12).Value = "12/06/2024 04:46:55"

# Row 232
$ws1.Cells.Item(232, 1).Value = 45441.55208333334
$ws1.Cells.Item(232, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(232, 2).Value = "12-06-2024 09:15:00"
$ws1.Cells.Item(232, 3).Value = "hour"
$ws1.Cells.Item(232, 4).Value = "TATASTEEL.NS"
$ws1.Cells.Item(232, 5).Value = 45415.38541666666
$ws1.Cells.Item(232, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(232, 6).Value = 170.75
$ws1.Cells.Item(232, 7).Value = 45436.38541666666
$ws1.Cells.Item(232, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(232, 8).Value = 177.5500030517578
$ws1.Cells.Item(232, 9).Value = 45440.38541666666
$ws1.Cells.Item(232, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(232, 10).Value = 177.5
$ws1.Cells.Item(232, 11).Value = "High"
$ws1.Cells.Item(232, 12).Value = "12/06/2024 04:46:55"

# Row 233
$ws1.Cells.Item(233, 1).Value = 45441.55208333334
$ws1.Cells.Item(233, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(233, 2).Value = "12-06-2024 09:15:00"
$ws1.Cells.Item(233, 3).Value = "hour"
$ws1.Cells.Item(233, 4).Value = "TATASTEEL.NS"
$ws1.Cells.Item(233, 5).Value = 45415.38541666666
$ws1.Cells.Item(233, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(233, 6).Value = 170.75
$ws1.Cells.Item(233, 7).Value = 45434.38541666666
$ws1.Cells.Item(233, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(233, 8).Value = 175.4499969482422
$ws1.Cells.Item(233, 9).Value = 45440.38541666666
$ws1.Cells.Item(233, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(233, 10).Value = 177.5
$ws1.Cells.Item(233, 11).Value = "High"
$ws1.Cells.Item(233, 12).Value = "12/06/2024 04:46:55"

# Row 234
$ws1.Cells.Item(234, 1).Value = 45450.55208333334
$ws1.Cells.Item(234, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(234, 2).Value = "12-06-2024 10:15:00"
$ws1.Cells.Item(234, 3).Value = "hour"
$ws1.Cells.Item(234, 4).Value = "FEDERALBNK.NS"
$ws1.Cells.Item(234, 5).Value = 45427.38541666666
$ws1.Cells.Item(234, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(234, 6).Value = 164.3500061035156
$ws1.Cells.Item(234, 7).Value = 45433.51041666666
$ws1.Cells.Item(234, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(234, 8).Value = 165.1999969482422
$ws1.Cells.Item(234, 9).Value = 45446.38541666666
$ws1.Cells.Item(234, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(234, 10).Value = 166.9499969482422
$ws1.Cells.Item(234, 11).Value = "High"
$ws1.Cells.Item(234, 12).Value = "12/06/2024 04:46:55"

# Row 235
$ws1.Cells.Item(235, 1).Value = 45450.63541666666
$ws1.Cells.Item(235, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(235, 2).Value = "12-06-2024 09:15:00"
$ws1.Cells.Item(235, 3).Value = "hour"
$ws1.Cells.Item(235, 4).Value = "PETRONET.NS"
$ws1.Cells.Item(235, 5).Value = 45415.38541666666
$ws1.Cells.Item(235, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(235, 6).Value = 322.75
$ws1.Cells.Item(235, 7).Value = 45429.38541666666
$ws1.Cells.Item(235, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(235, 8).Value = 320.1499938964844
$ws1.Cells.Item(235, 9).Value = 45446.59375
$ws1.Cells.Item(235, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(235, 10).Value = 317.7000122070312
$ws1.Cells.Item(235, 11).Value = "High"
$ws1.Cells.Item(235, 12).Value = "12/06/2024 04:46:55"

# Row 236
$ws1.Cells.Item(236, 1).Value = 45440.59375
$ws1.Cells.Item(236, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(236, 2).Value = "12-06-2024 09:15:00"
$ws1.Cells.Item(236, 3).Value = "hour"
$ws1.Cells.Item(236, 4).Value = "RAIN.NS"
$ws1.Cells.Item(236, 5).Value = 45408.46875
$ws1.Cells.Item(236, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(236, 6).Value = 184.5
$ws1.Cells.Item(236, 7).Value = 45434.42708333334
$ws1.Cells.Item(236, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(236, 8).Value = 174.75
$ws1.Cells.Item(236, 9).Value = 45439.42708333334
$ws1.Cells.Item(236, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(236, 10).Value = 173.6499938964844
$ws1.Cells.Item(236, 11).Value = "High"
$ws1.Cells.Item(236, 12).Value = "12/06/2024 04:46:55"

# Row 237
$ws1.Cells.Item(237, 1).Value = 45434.55208333334
$ws1.Cells.Item(237, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(237, 2).Value = "12-06-2024 09:15:00"
$ws1.Cells.Item(237, 3).Value = "hour"
$ws1.Cells.Item(237, 4).Value = "NECLIFE.NS"
$ws1.Cells.Item(237, 5).Value = 45415.38541666666
$ws1.Cells.Item(237, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(237, 6).Value = 35.90000152587891
$ws1.Cells.Item(237, 7).Value = 45418.59375
$ws1.Cells.Item(237, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(237, 8).Value = 35.75
$ws1.Cells.Item(237, 9).Value = 45433.38541666666
$ws1.Cells.Item(237, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(237, 10).Value = 34.90000152587891
$ws1.Cells.Item(237, 11).Value = "High"
$ws1.Cells.Item(237, 12).Value = "12/06/2024 04:46:55"

# Row 238
$ws1.Cells.Item(238, 1).Value = 45447.55208333334
$ws1.Cells.Item(238, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(238, 2).Value = "12-06-2024 09:15:00"
$ws1.Cells.Item(238, 3).Value = "hour"
$ws1.Cells.Item(238, 4).Value = "MCX.NS"
$ws1.Cells.Item(238, 5).Value = 45427.38541666666
$ws1.Cells.Item(238, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(238, 6).Value = 3963.39990234375
$ws1.Cells.Item(238, 7).Value = 45440.38541666666
$ws1.Cells.Item(238, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(238, 8).Value = 3882.949951171875
$ws1.Cells.Item(238, 9).Value = 45446.38541666666
$ws1.Cells.Item(238, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(238, 10).Value = 3808.14990234375
$ws1.Cells.Item(238, 11).Value = "High"
$ws1.Cells.Item(238, 12).Value = "12/06/2024 04:46:55"

# Row 239
$ws1.Cells.Item(239, 1).Value = 45447.59375
$ws1.Cells.Item(239, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(239, 2).Value = "12-06-2024 09:15:00"
$ws1.Cells.Item(239, 3).Value = "hour"
$ws1.Cells.Item(239, 4).Value = "FAZE3Q.NS"
$ws1.Cells.Item(239, 5).Value = 45434.38541666666
$ws1.Cells.Item(239, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(239, 6).Value = 440.3500061035156
$ws1.Cells.Item(239, 7).Value = 45435.55208333334
$ws1.Cells.Item(239, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(239, 8).Value = 438.5
$ws1.Cells.Item(239, 9).Value = 45435.59375
$ws1.Cells.Item(239, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(239, 10).Value = 438.5
$ws1.Cells.Item(239, 11).Value = "High"
$ws1.Cells.Item(239, 12).Value = "12/06/2024 04:46:55"

# Sheet 2 new rows
# Row 56
$ws2.Cells.Item(56, 1).Value = 45448.59375
$ws2.Cells.Item(56, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(56, 2).Value = "12-06-2024 09:15:00"
$ws2.Cells.Item(56, 3).Value = "hour"
$ws2.Cells.Item(56, 4).Value = "SHILCTECH.BO"
$ws2.Cells.Item(56, 5).Value = 45447.38541666666
$ws2.Cells.Item(56, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(56, 6).Value = 5350
$ws2.Cells.Item(56, 7).Value = 45447.42708333334
$ws2.Cells.Item(56, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(56, 8).Value = 5350
$ws2.Cells.Item(56, 9).Value = "High"
$ws2.Cells.Item(56, 10).Value = "12/06/2024 04:46:55"

# Row 57
$ws2.Cells.Item(57, 1).Value = 45454.55208333334
$ws2.Cells.Item(57, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(57, 2).Value = "12-06-2024 09:15:00"
$ws2.Cells.Item(57, 3).Value = "hour"
$ws2.Cells.Item(57, 4).Value = "KPEL.BO"
$ws2.Cells.Item(57, 5).Value = 45442.38541666666
$ws2.Cells.Item(57, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(57, 6).Value = 437.5
$ws2.Cells.Item(57, 7).Value = 45453.38541666666
$ws2.Cells.Item(57, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(57, 8).Value = 438
$ws2.Cells.Item(57, 9).Value = "High"
$ws2.Cells.Item(57, 10).Value = "12/06/2024 04:46:55"

# Row 58
$ws2.Cells.Item(58, 1).Value = 45454.55208333334
$ws2.Cells.Item(58, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(58, 2).Value = "12-06-2024 09:15:00"
$ws2.Cells.Item(58, 3).Value = "hour"
$ws2.Cells.Item(58, 4).Value = "KPEL.BO"
$ws2.Cells.Item(58, 5).Value = 45446.38541666666
$ws2.Cells.Item(58, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(58, 6).Value = 436.2000122070312
$ws2.Cells.Item(58, 7).Value = 45453.38541666666
$ws2.Cells.Item(58, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(58, 8).Value = 438
$ws2.Cells.Item(58, 9).Value = "High"
$ws2.Cells.Item(58, 10).Value = "12/06/2024 04:46:55"

# Row 59
$ws2.Cells.Item(59, 1).Value = 45446.38541666666
$ws2.Cells.Item(59, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(59, 2).Value = "12-06-2024 10:15:00"
$ws2.Cells.Item(59, 3).Value = "hour"
$ws2.Cells.Item(59, 4).Value = "SHANKARA.NS"
$ws2.Cells.Item(59, 5).Value = 45433.46875
$ws2.Cells.Item(59, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(59, 6).Value = 715
$ws2.Cells.Item(59, 7).Value = 45442.51041666666
$ws2.Cells.Item(59, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(59, 8).Value = 712.5
$ws2.Cells.Item(59, 9).Value = "High"
$ws2.Cells.Item(59, 10).Value = "12/06/2024 04:46:55"

# Row 60
$ws2.Cells.Item(60, 1).Value = 45447.55208333334
$ws2.Cells.Item(60, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(60, 2).Value = "12-06-2024 10:15:00"
$ws2.Cells.Item(60, 3).Value = "hour"
$ws2.Cells.Item(60, 4).Value = "SHANKARA.NS"
$ws2.Cells.Item(60, 5).Value = 45433.46875
$ws2.Cells.Item(60, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(60, 6).Value = 715
$ws2.Cells.Item(60, 7).Value = 45446.38541666666
$ws2.Cells.Item(60, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(60, 8).Value = 711.7999877929688
$ws2.Cells.Item(60, 9).Value = "High"
$ws2.Cells.Item(60, 10).Value = "12/06/2024 04:46:55"

# Sheet 3 new rows
# Row 758
$ws3.Cells.Item(758, 1).Value = "TRIL.BO"
$ws3.Cells.Item(758, 2).Value = 45447.38541666666
$ws3.Cells.Item(758, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(758, 3).Value = 771.2999877929688
$ws3.Cells.Item(758, 4).Value = 706.5499877929688
$ws3.Cells.Item(758, 5).Value = 731
$ws3.Cells.Item(758, 6).Value = "High"
$ws3.Cells.Item(758, 7).Value = 771.2999877929688
$ws3.Cells.Item(758, 8).Value = "hour"
$ws3.Cells.Item(758, 9).Value = "12-06-2024 09:15:00"
$ws3.Cells.Item(758, 10).Value = 807.5999755859375
$ws3.Cells.Item(758, 11).Value = 769.1500244140625
$ws3.Cells.Item(758, 12).Value = "12/06/2024 04:46:55"

# Row 759
$ws3.Cells.Item(759, 1).Value = "SHILCTECH.BO"
$ws3.Cells.Item(759, 2).Value = 45447.38541666666
$ws3.Cells.Item(759, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(759, 3).Value = 5350
$ws3.Cells.Item(759, 4).Value = 5029.14990234375
$ws3.Cells.Item(759, 5).Value = 5200
$ws3.Cells.Item(759, 6).Value = "High"
$ws3.Cells.Item(759, 7).Value = 5350
$ws3.Cells.Item(759, 8).Value = "hour"
$ws3.Cells.Item(759, 9).Value = "12-06-2024 09:15:00"
$ws3.Cells.Item(759, 10).Value = 5400
$ws3.Cells.Item(759, 11).Value = 5290
$ws3.Cells.Item(759, 12).Value = "12/06/2024 04:46:55"

# Row 760
$ws3.Cells.Item(760, 1).Value = "SHILCTECH.BO"
$ws3.Cells.Item(760, 2).Value = 45447.42708333334
$ws3.Cells.Item(760, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(760, 3).Value = 5350
$ws3.Cells.Item(760, 4).Value = 5029.14990234375
$ws3.Cells.Item(760, 5).Value = 5029.14990234375
$ws3.Cells.Item(760, 6).Value = "High"
$ws3.Cells.Item(760, 7).Value = 5350
$ws3.Cells.Item(760, 8).Value = "hour"
$ws3.Cells.Item(760, 9).Value = "12-06-2024 09:15:00"
$ws3.Cells.Item(760, 10).Value = 5400
$ws3.Cells.Item(760, 11).Value = 5290
$ws3.Cells.Item(760, 12).Value = "12/06/2024 04:46:55"

# Row 761
$ws3.Cells.Item(761, 1).Value = "SHILCTECH.BO"
$ws3.Cells.Item(761, 2).Value = 45449.38541666666
$ws3.Cells.Item(761, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(761, 3).Value = 5394
$ws3.Cells.Item(761, 4).Value = 5100
$ws3.Cells.Item(761, 5).Value = 5200
$ws3.Cells.Item(761, 6).Value = "High"
$ws3.Cells.Item(761, 7).Value = 5394
$ws3.Cells.Item(761, 8).Value = "hour"
$ws3.Cells.Item(761, 9).Value = "12-06-2024 09:15:00"
$ws3.Cells.Item(761, 10).Value = 5400
$ws3.Cells.Item(761, 11).Value = 5290
$ws3.Cells.Item(761, 12).Value = "12/06/2024 04:46:55"

# Row 762
$ws3.Cells.Item(762, 1).Value = "SPELS.BO"
$ws3.Cells.Item(762, 2).Value = 45434.42708333334
$ws3.Cells.Item(762, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(762, 3).Value = 141.6999969482422
$ws3.Cells.Item(762, 4).Value = 141.6999969482422
$ws3.Cells.Item(762, 5).Value = 141.6999969482422
$ws3.Cells.Item(762, 6).Value = "High"
$ws3.Cells.Item(762, 7).Value = 141.6999969482422
$ws3.Cells.Item(762, 8).Value = "hour"
$ws3.Cells.Item(762, 9).Value = "12-06-2024 09:15:00"
$ws3.Cells.Item(762, 10).Value = 143.9499969482422
$ws3.Cells.Item(762, 11).Value = 139.9499969482422
$ws3.Cells.Item(762, 12).Value = "12/06/2024 04:46:55"

# Row 763
$ws3.Cells.Item(763, 1).Value = "SPELS.BO"
$ws3.Cells.Item(763, 2).Value = 45434.46875
$ws3.Cells.Item(763, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(763, 3).Value = 141.6999969482422
$ws3.Cells.Item(763, 4).Value = 141.6999969482422
$ws3.Cells.Item(763, 5).Value = 141.6999969482422
$ws3.Cells.Item(763, 6).Value = "High"
$ws3.Cells.Item(763, 7).Value = 141.6999969482422
$ws3.Cells.Item(763, 8).Value = "hour"
$ws3.Cells.Item(763, 9).Value = "12-06-2024 09:15:00"
$ws3.Cells.Item(763, 10).Value = 143.9499969482422
$ws3.Cells.Item(763, 11).Value = 139.9499969482422
$ws3.Cells.Item(763, 12).Value = "12/06/2024 04:46:55"

# Row 764
$ws3.Cells.Item(764, 1).Value = "SIKA.BO"
$ws3.Cells.Item(764, 2).Value = 45436.38541666666
$ws3.Cells.Item(764, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(764, 3).Value = 2749
$ws3.Cells.Item(764, 4).Value = 2680
$ws3.Cells.Item(764, 5).Value = 2707
$ws3.Cells.Item(764, 6).Value = "High"
$ws3.Cells.Item(764, 7).Value = 2749
$ws3.Cells.Item(764, 8).Value = "hour"
$ws3.Cells.Item(764, 9).Value = "12-06-2024 09:15:00"
$ws3.Cells.Item(764, 10).Value = 2792
$ws3.Cells.Item(764, 11).Value = 2690
$ws3.Cells.Item(764, 12).Value = "12/06/2024 04:46:55"

# Row 765
$ws3.Cells.Item(765, 1).Value = "SIL.NS"
$ws3.Cells.Item(765, 2).Value = 45446.38541666666
$ws3.Cells.Item(765, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(765, 3).Value = 24.89999961853027
$ws3.Cells.Item(765, 4).Value = 23.45000076293945
$ws3.Cells.Item(765, 5).Value = 23.85000038146973
$ws3.Cells.Item(765, 6).Value = "High"
$ws3.Cells.Item(765, 7).Value = 24.89999961853027
$ws3.Cells.Item(765, 8).Value = "hour"
$ws3.Cells.Item(765, 9).Value = "12-06-2024 09:15:00"
$ws3.Cells.Item(765, 10).Value = 25.35000038146973
$ws3.Cells.Item(765, 11).Value = 24.6200008392334
$ws3.Cells.Item(765, 12).Value = "12/06/2024 04:46:55"

# Row 766
$ws3.Cells.Item(766, 1).Value = "KPEL.BO"
$ws3.Cells.Item(766, 2).Value = 45442.38541666666
$ws3.Cells.Item(766, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(766, 3).Value = 437.5
$ws3.Cells.Item(766, 4).Value = 420.2000122070312
$ws3.Cells.Item(766, 5).Value = 430.25
$ws3.Cells.Item(766, 6).Value = "High"
$ws3.Cells.Item(766, 7).Value = 437.5
$ws3.Cells.Item(766, 8).Value = "hour"
$ws3.Cells.Item(766, 9).Value = "12-06-2024 09:15:00"
$ws3.Cells.Item(766, 10).Value = 445
$ws3.Cells.Item(766, 11).Value = 436.8999938964844
$ws3.Cells.Item(766, 12).Value = "12/06/2024 04:46:55"

# Row 767
$ws3.Cells.Item(767, 1).Value = "KPEL.BO"
$ws3.Cells.Item(767, 2).Value = 45453.38541666666
$ws3.Cells.Item(767, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(767, 3).Value = 438
$ws3.Cells.Item(767, 4).Value = 420.1000061035156
$ws3.Cells.Item(767, 5).Value = 424.3500061035156
$ws3.Cells.Item(767, 6).Value = "High"
$ws3.Cells.Item(767, 7).Value = 438
$ws3.Cells.Item(767, 8).Value = "hour"
$ws3.Cells.Item(767, 9).Value = "12-06-2024 09:15:00"
$ws3.Cells.Item(767, 10).Value = 445
$ws3.Cells.Item(767, 11).Value = 436.8999938964844
$ws3.Cells.Item(767, 12).Value = "12/06/2024 04:46:55"

# Row 768
$ws3.Cells.Item(768, 1).Value = "UDS.NS"
$ws3.Cells.Item(768, 2).Value = 45439.38541666666
$ws3.Cells.Item(768, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(768, 3).Value = 303
$ws3.Cells.Item(768, 4).Value = 295.2999877929688
$ws3.Cells.Item(768, 5).Value = 297.5499877929688
$ws3.Cells.Item(768, 6).Value = "High"
$ws3.Cells.Item(768, 7).Value = 303
$ws3.Cells.Item(768, 8).Value = "hour"
$ws3.Cells.Item(768, 9).Value = "12-06-2024 09:15:00"
$ws3.Cells.Item(768, 10).Value = 303.3999938964844
$ws3.Cells.Item(768, 11).Value = 297.8500061035156
$ws3.Cells.Item(768, 12).Value = "12/06/2024 04:46:55"

# Row 769
$ws3.Cells.Item(769, 1).Value = "UNIAUTO.BO"
$ws3.Cells.Item(769, 2).Value = 45433.38541666666
$ws3.Cells.Item(769, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(769, 3).Value = 168.3999938964844
$ws3.Cells.Item(769, 4).Value = 163
$ws3.Cells.Item(769, 5).Value = 163
$ws3.Cells.Item(769, 6).Value = "High"
$ws3.Cells.Item(769, 7).Value = 168.3999938964844
$ws3.Cells.Item(769, 8).Value = "hour"
$ws3.Cells.Item(769, 9).Value = "12-06-2024 09:15:00"
$ws3.Cells.Item(769, 10).Value = 169.6999969482422
$ws3.Cells.Item(769, 11).Value = 167
$ws3.Cells.Item(769, 12).Value = "12/06/2024 04:46:55"
